# Applies the "Add files via upload" edit to NPS-Presentacion.pptx:
#  - refresh the four datetimeFigureOut placeholders (slide master, the
#    "Title and Content" layout, notes master, handout master) from
#    7/24/2021 to 7/27/2021
#  - retitle the title slide ("Welcome to PowerPoint" -> "NPS")
#  - retitle the title slide subtitle ("5 tips for a simpler way to work"
#    -> "Julio 2021 Datasphere", typed as three runs so the "Julio"
#    portion keeps its own formatting run like the captured edit)

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($Shapes, [string]$NewText)

    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $sh = $Shapes.Item($i)
        $phType = $null
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        # ppPlaceholderDate = 16
        if ($phType -eq 16 -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $NewText
            return $true
        }
    }
    return $false
}

$newDate = "7/27/2021"

# Slide master "Date Placeholder"
Set-DatePlaceholderText -Shapes $p.SlideMaster.Shapes -NewText $newDate | Out-Null

# Every slide layout off the master that carries its own date placeholder
# (only "Title and Content" has one in this deck, but loop defensively)
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText -Shapes $layouts.Item($li).Shapes -NewText $newDate | Out-Null
}

# Notes master + handout master date placeholders
Set-DatePlaceholderText -Shapes $p.NotesMaster.Shapes -NewText $newDate | Out-Null
Set-DatePlaceholderText -Shapes $p.HandoutMaster.Shapes -NewText $newDate | Out-Null

# Slide 1 (title slide): update title + subtitle text
$slide1 = $p.Slides.Item(1)

$titleShape = $null
$subtitleShape = $null
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $sh = $slide1.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    $phType = $null
    try { $phType = $sh.PlaceholderFormat.Type } catch {}
    if ($phType -eq 13 -or $phType -eq 1) { $titleShape = $sh }       # ctrTitle / title
    elseif ($phType -eq 4) { $subtitleShape = $sh }                   # subtitle
}
if (-not $titleShape) { $titleShape = $slide1.Shapes.Item(1) }
if (-not $subtitleShape) { $subtitleShape = $slide1.Shapes.Item(2) }

$titleShape.TextFrame.TextRange.Text = "NPS"

$subtitleRange = $subtitleShape.TextFrame.TextRange
$subtitleRange.Text = "J"
$subtitleRange.InsertAfter("ulio")
$subtitleRange.InsertAfter(" 2021 Datasphere")

Write-Output "Edit applied"
